{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the copyright/footer line right after it, and the blank paragraph that\n// separates them from the preceding \"Requisitos\" list - this matches the\n// Jekyll site rebuild that dropped the page-chrome paragraphs scraped along\n// with the real Jupiter-catalog content, while leaving the blank spacer\n// paragraph that sits right before the trailing page-break paragraph intact.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two text paragraphs that must be removed.\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (jupiterIndex === -1 && text.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n  }\n  if (copyrightIndex === -1 && text.indexOf(\"Powered by Jekyll\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (jupiterIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\"Could not locate the 'Ver no Jupiter' / copyright paragraphs to remove.\");\n}\n\n// The blank paragraph immediately preceding \"Ver no Jupiter\" is the\n// separator that gets removed together with the two text paragraphs.\nlet blankIndex = jupiterIndex - 1;\nif (blankIndex < 0 || items[blankIndex].text !== \"\") {\n  blankIndex = -1;\n}\n\n// Delete from the end of the run backwards so earlier indices stay valid.\nitems[copyrightIndex].delete();\nitems[jupiterIndex].delete();\nif (blankIndex !== -1) {\n  items[blankIndex].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, the\n# copyright/footer line right after it, and the blank separator paragraph\n# right before them - the Jekyll site rebuild dropped these page-chrome\n# paragraphs while keeping the blank spacer paragraph that precedes the\n# trailing page-break paragraph.\n$d = $word.ActiveDocument\n\n$jupiterIndex = -1\n$copyrightIndex = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIndex -eq -1 -and $text -like \"*Ver no Jupiter*\") {\n        $jupiterIndex = $i\n    }\n    if ($copyrightIndex -eq -1 -and $text -like \"*Powered by Jekyll*\") {\n        $copyrightIndex = $i\n    }\n}\n\nif ($jupiterIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"Could not locate the 'Ver no Jupiter' / copyright paragraphs to remove.\"\n}\n\n$blankIndex = $jupiterIndex - 1\nif ($blankIndex -ge 1) {\n    $blankText = $d.Paragraphs.Item($blankIndex).Range.Text\n    if ($blankText -ne \"`r\") {\n        $blankIndex = -1\n    }\n} else {\n    $blankIndex = -1\n}\n\n# Delete back-to-front so earlier indices stay valid.\n$d.Paragraphs.Item($copyrightIndex).Range.Delete()\n$d.Paragraphs.Item($jupiterIndex).Range.Delete()\nif ($blankIndex -ne -1) {\n    $d.Paragraphs.Item($blankIndex).Range.Delete()\n}\n"}
